$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.022639180589272
$ws.Cells.Item(2, 4).Value = 1.025283847741079
$ws.Cells.Item(2, 5).Value = 1.023349496787121
$ws.Cells.Item(2, 6).Value = 1.021090700653377
$ws.Cells.Item(2, 9).Value = 1.028772394481805
$ws.Cells.Item(2, 10).Value = 1.027823887186856
$ws.Cells.Item(2, 11).Value = 1.028109870440689
$ws.Cells.Item(2, 12).Value = 1.026181195488365
$ws.Cells.Item(2, 13).Value = 1.023929057002742
$ws.Cells.Item(2, 14).Value = 1.029283514999711
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.023733011756376
$ws.Cells.Item(3, 4).Value = 1.026290199793475
$ws.Cells.Item(3, 5).Value = 1.024281688079975
$ws.Cells.Item(3, 6).Value = 1.022824625065921
$ws.Cells.Item(3, 9).Value = 1.02890981021903
$ws.Cells.Item(3, 10).Value = 1.028555104608814
$ws.Cells.Item(3, 11).Value = 1.028922952208558
$ws.Cells.Item(3, 12).Value = 1.026919905077879
$ws.Cells.Item(3, 13).Value = 1.025466820857279
$ws.Cells.Item(3, 14).Value = 1.030015770834279
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.024440262294551
$ws.Cells.Item(4, 4).Value = 1.026941148960525
$ws.Cells.Item(4, 5).Value = 1.024884785732375
$ws.Cells.Item(4, 6).Value = 1.02394585655939
$ws.Cells.Item(4, 9).Value = 1.028996988339065
$ws.Cells.Item(4, 10).Value = 1.0290272015447
$ws.Cells.Item(4, 11).Value = 1.029448240513462
$ws.Cells.Item(4, 12).Value = 1.027397188519306
$ws.Cells.Item(4, 13).Value = 1.026460692095375
$ws.Cells.Item(4, 14).Value = 1.030488538201964
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.024737465237574
$ws.Cells.Item(5, 4).Value = 1.027214754807147
$ws.Cells.Item(5, 5).Value = 1.025138306886164
$ws.Cells.Item(5, 6).Value = 1.02441705537307
$ws.Cells.Item(5, 9).Value = 1.029033221845327
$ws.Cells.Item(5, 10).Value = 1.029225421133059
$ws.Cells.Item(5, 11).Value = 1.029668874357661
$ws.Cells.Item(5, 12).Value = 1.027597669344476
$ws.Cells.Item(5, 13).Value = 1.026878244308405
$ws.Cells.Item(5, 14).Value = 1.030687039284874
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.024787359637374
$ws.Cells.Item(6, 4).Value = 1.027260691353823
$ws.Cells.Item(6, 5).Value = 1.025180872995954
$ws.Cells.Item(6, 6).Value = 1.024496162192127
$ws.Cells.Item(6, 9).Value = 1.029039281210204
$ws.Cells.Item(6, 10).Value = 1.029258688454325
$ws.Cells.Item(6, 11).Value = 1.029705908216325
$ws.Cells.Item(6, 12).Value = 1.027631321070729
$ws.Cells.Item(6, 13).Value = 1.026948337476009
$ws.Cells.Item(6, 14).Value = 1.030720353849552
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.024444234024201
$ws.Cells.Item(7, 4).Value = 1.026944805102525
$ws.Cells.Item(7, 5).Value = 1.024888173376087
$ws.Cells.Item(7, 6).Value = 1.023952153386058
$ws.Cells.Item(7, 9).Value = 1.02899747412805
$ws.Cells.Item(7, 10).Value = 1.029029851145113
$ws.Cells.Item(7, 11).Value = 1.029451189407539
$ws.Cells.Item(7, 12).Value = 1.027399868017147
$ws.Cells.Item(7, 13).Value = 1.026466272504938
$ws.Cells.Item(7, 14).Value = 1.030491191565113
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.023008957201222
$ws.Cells.Item(8, 4).Value = 1.025623997151068
$ws.Cells.Item(8, 5).Value = 1.023664554966696
$ws.Cells.Item(8, 6).Value = 1.021676845662572
$ws.Cells.Item(8, 9).Value = 1.028819194951577
$ws.Cells.Item(8, 10).Value = 1.028071223669672
$ws.Cells.Item(8, 11).Value = 1.028384827889886
$ws.Cells.Item(8, 12).Value = 1.026430993745137
$ws.Cells.Item(8, 13).Value = 1.02444899726129
$ws.Cells.Item(8, 14).Value = 1.0295312027287
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.020475654810981
$ws.Cells.Item(9, 4).Value = 1.023294756682117
$ws.Cells.Item(9, 5).Value = 1.021507627738892
$ws.Cells.Item(9, 6).Value = 1.01766145411696
$ws.Cells.Item(9, 9).Value = 1.028491717727952
$ws.Cells.Item(9, 10).Value = 1.026373900722358
$ws.Cells.Item(9, 11).Value = 1.026499345363279
$ws.Cells.Item(9, 12).Value = 1.024718217757174
$ws.Cells.Item(9, 13).Value = 1.020885038486687
$ws.Cells.Item(9, 14).Value = 1.027831469388115
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.01878385780994
$ws.Cells.Item(10, 4).Value = 1.021740623306226
$ws.Cells.Item(10, 5).Value = 1.020069099333576
$ws.Cells.Item(10, 6).Value = 1.01497993940633
$ws.Cells.Item(10, 9).Value = 1.028264427670716
$ws.Cells.Item(10, 10).Value = 1.025236820518306
$ws.Cells.Item(10, 11).Value = 1.025237954733939
$ws.Cells.Item(10, 12).Value = 1.023572599750691
$ws.Cells.Item(10, 13).Value = 1.01850236053351
$ws.Cells.Item(10, 14).Value = 1.026692774399749
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.018050566339804
$ws.Cells.Item(11, 4).Value = 1.021067335376068
$ws.Cells.Item(11, 5).Value = 1.019446047473009
$ws.Cells.Item(11, 6).Value = 1.013817603893067
$ws.Cells.Item(11, 9).Value = 1.028163877442753
$ws.Cells.Item(11, 10).Value = 1.024743120209659
$ws.Cells.Item(11, 11).Value = 1.024690694145812
$ws.Cells.Item(11, 12).Value = 1.023075624041452
$ws.Cells.Item(11, 13).Value = 1.017468937456996
$ws.Cells.Item(11, 14).Value = 1.026198372980036
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.017778076538747
$ws.Cells.Item(12, 4).Value = 1.020817193561582
$ws.Cells.Item(12, 5).Value = 1.019214592786641
$ws.Cells.Item(12, 6).Value = 1.013385666578458
$ws.Cells.Item(12, 9).Value = 1.028126207948189
$ws.Cells.Item(12, 10).Value = 1.02455953523181
$ws.Cells.Item(12, 11).Value = 1.024487254568623
$ws.Cells.Item(12, 12).Value = 1.022890886034403
$ws.Cells.Item(12, 13).Value = 1.017084812273505
$ws.Cells.Item(12, 14).Value = 1.026014527290461
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.017836531669408
$ws.Cells.Item(13, 4).Value = 1.020870852245863
$ws.Cells.Item(13, 5).Value = 1.019264241780214
$ws.Cells.Item(13, 6).Value = 1.013478327554889
$ws.Cells.Item(13, 9).Value = 1.028134302698747
$ws.Cells.Item(13, 10).Value = 1.02459892403483
$ws.Cells.Item(13, 11).Value = 1.02453090043775
$ws.Cells.Item(13, 12).Value = 1.022930519283682
$ws.Cells.Item(13, 13).Value = 1.017167220608468
$ws.Cells.Item(13, 14).Value = 1.026053972030099
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.01802804456253
$ws.Cells.Item(14, 4).Value = 1.021046659661933
$ws.Cells.Item(14, 5).Value = 1.019426915885102
$ws.Cells.Item(14, 6).Value = 1.013781903841921
$ws.Cells.Item(14, 9).Value = 1.028160770213962
$ws.Cells.Item(14, 10).Value = 1.02472794916376
$ws.Cells.Item(14, 11).Value = 1.024673881095254
$ws.Cells.Item(14, 12).Value = 1.023060356386798
$ws.Cells.Item(14, 13).Value = 1.017437191026346
$ws.Cells.Item(14, 14).Value = 1.026183180389511
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.018146026962035
$ws.Cells.Item(15, 4).Value = 1.02115497337213
$ws.Cells.Item(15, 5).Value = 1.019527141347664
$ws.Cells.Item(15, 6).Value = 1.013968921167281
$ws.Cells.Item(15, 9).Value = 1.028177035222579
$ws.Cells.Item(15, 10).Value = 1.024807418889414
$ws.Cells.Item(15, 11).Value = 1.024761954591722
$ws.Cells.Item(15, 12).Value = 1.023140334832348
$ws.Cells.Item(15, 13).Value = 1.017603493206482
$ws.Cells.Item(15, 14).Value = 1.026262762971291
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.018832508293598
$ws.Cells.Item(16, 4).Value = 1.021785299924812
$ws.Cells.Item(16, 5).Value = 1.020110445743154
$ws.Cells.Item(16, 6).Value = 1.015057053140778
$ws.Cells.Item(16, 9).Value = 1.028271055908233
$ws.Cells.Item(16, 10).Value = 1.025269557456174
$ws.Cells.Item(16, 11).Value = 1.025274251912732
$ws.Cells.Item(16, 12).Value = 1.023605562970262
$ws.Cells.Item(16, 13).Value = 1.01857090864904
$ws.Cells.Item(16, 14).Value = 1.026725557827824
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.019262921968343
$ws.Cells.Item(17, 4).Value = 1.022180595290063
$ws.Cells.Item(17, 5).Value = 1.020476293397657
$ws.Cells.Item(17, 6).Value = 1.015739274989882
$ws.Cells.Item(17, 9).Value = 1.028329461533777
$ws.Cells.Item(17, 10).Value = 1.02555908551611
$ws.Cells.Item(17, 11).Value = 1.025595314681841
$ws.Cells.Item(17, 12).Value = 1.023897142081668
$ws.Cells.Item(17, 13).Value = 1.019177279050701
$ws.Cells.Item(17, 14).Value = 1.027015497050815
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.019513904519793
$ws.Cells.Item(18, 4).Value = 1.022411131916662
$ws.Cells.Item(18, 5).Value = 1.020689670852969
$ws.Cells.Item(18, 6).Value = 1.016137086269685
$ws.Cells.Item(18, 9).Value = 1.02836332285039
$ws.Cells.Item(18, 10).Value = 1.025727833451318
$ws.Cells.Item(18, 11).Value = 1.025782482039068
$ws.Cells.Item(18, 12).Value = 1.024067127013121
$ws.Cells.Item(18, 13).Value = 1.019530800498329
$ws.Cells.Item(18, 14).Value = 1.027184484627446
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.019599471212865
$ws.Cells.Item(19, 4).Value = 1.022489733498542
$ws.Cells.Item(19, 5).Value = 1.020762424502284
$ws.Cells.Item(19, 6).Value = 1.016272710077979
$ws.Cells.Item(19, 9).Value = 1.028374833804428
$ws.Cells.Item(19, 10).Value = 1.025785350329929
$ws.Cells.Item(19, 11).Value = 1.025846283834855
$ws.Cells.Item(19, 12).Value = 1.024125072562681
$ws.Cells.Item(19, 13).Value = 1.019651314612729
$ws.Cells.Item(19, 14).Value = 1.027242083186621
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.019216749963369
$ws.Cells.Item(20, 4).Value = 1.022138187212083
$ws.Cells.Item(20, 5).Value = 1.020437043000187
$ws.Cells.Item(20, 6).Value = 1.015666091241382
$ws.Cells.Item(20, 9).Value = 1.028323216441769
$ws.Cells.Item(20, 10).Value = 1.025528035240382
$ws.Cells.Item(20, 11).Value = 1.025560878366108
$ws.Cells.Item(20, 12).Value = 1.02386586754063
$ws.Cells.Item(20, 13).Value = 1.019112238292942
$ws.Cells.Item(20, 14).Value = 1.026984402680134
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.017971651921623
$ws.Cells.Item(21, 4).Value = 1.020994890230138
$ws.Cells.Item(21, 5).Value = 1.019379013131464
$ws.Cells.Item(21, 6).Value = 1.013692513636592
$ws.Cells.Item(21, 9).Value = 1.028152985044995
$ws.Cells.Item(21, 10).Value = 1.024689960089411
$ws.Cells.Item(21, 11).Value = 1.02463178135963
$ws.Cells.Item(21, 12).Value = 1.023022126448424
$ws.Cells.Item(21, 13).Value = 1.01735769888028
$ws.Cells.Item(21, 14).Value = 1.02614513736632
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.017188155973447
$ws.Cells.Item(22, 4).Value = 1.020275748169291
$ws.Cells.Item(22, 5).Value = 1.018713639556003
$ws.Cells.Item(22, 6).Value = 1.012450518478137
$ws.Cells.Item(22, 9).Value = 1.028044098382739
$ws.Cells.Item(22, 10).Value = 1.024161855297698
$ws.Cells.Item(22, 11).Value = 1.024046679348813
$ws.Cells.Item(22, 12).Value = 1.022490827362008
$ws.Cells.Item(22, 13).Value = 1.016253007616735
$ws.Cells.Item(22, 14).Value = 1.025616282605228
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.017603564837693
$ws.Cells.Item(23, 4).Value = 1.020657008624404
$ws.Cells.Item(23, 5).Value = 1.019066381234283
$ws.Cells.Item(23, 6).Value = 1.013109034133607
$ws.Cells.Item(23, 9).Value = 1.028101997248579
$ws.Cells.Item(23, 10).Value = 1.024441925503268
$ws.Cells.Item(23, 11).Value = 1.024356942894879
$ws.Cells.Item(23, 12).Value = 1.02277255594051
$ws.Cells.Item(23, 13).Value = 1.016838774694291
$ws.Cells.Item(23, 14).Value = 1.025896750542614
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.019237613322636
$ws.Cells.Item(24, 4).Value = 1.022157349697435
$ws.Cells.Item(24, 5).Value = 1.020454778612583
$ws.Cells.Item(24, 6).Value = 1.015699160188335
$ws.Cells.Item(24, 9).Value = 1.028326038965625
$ws.Cells.Item(24, 10).Value = 1.02554206592164
$ws.Cells.Item(24, 11).Value = 1.025576438973513
$ws.Cells.Item(24, 12).Value = 1.023879999432148
$ws.Cells.Item(24, 13).Value = 1.01914162791813
$ws.Cells.Item(24, 14).Value = 1.026998453286568
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.021131080750849
$ws.Cells.Item(25, 4).Value = 1.0238971447205
$ws.Cells.Item(25, 5).Value = 1.022065341686242
$ws.Cells.Item(25, 6).Value = 1.018700297272295
$ws.Cells.Item(25, 9).Value = 1.028577958905247
$ws.Cells.Item(25, 10).Value = 1.026813667978081
$ws.Cells.Item(25, 11).Value = 1.026987557237068
$ws.Cells.Item(25, 12).Value = 1.025161670068929
$ws.Cells.Item(25, 13).Value = 1.021807553916295
$ws.Cells.Item(25, 14).Value = 1.028271861163784
